$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = 2.12
$ws.Range("H5").Value = 3
$ws.Range("J5").Value = 1.42
$ws.Range("K5").Value = 2.47
$ws.Range("L5").Value = 2.2
$ws.Range("M5").Value = 1.53
$ws.Range("N5").Value = 1.47
$ws.Range("O5").Value = 2.35
$ws.Range("P5").Value = 1.91
$ws.Range("Q5").Value = 1.7
$ws.Range("R5").Value = 6.1
$ws.Range("S5").Value = 9.25
$ws.Range("T5").Value = 9
$ws.Range("V5").Value = 20
$ws.Range("W5").Value = 37
$ws.Range("X5").Value = 7.2
$ws.Range("Y5").Value = 5.9
$ws.Range("Z5").Value = 16.5
$ws.Range("AA5").Value = 100
$ws.Range("AB5").Value = 8.5
$ws.Range("AG5").Value = 50
$ws.Range("AH5").Value = 900
